# Apply the "Updated symbol list" data refresh to the cryptos sheet.
# The Price column (D) holds numeric-looking values that were authored as
# literal text (inline strings), so each numeric update first forces the
# cell's number format to Text ("@") before assigning the string - this
# mirrors how Excel COM automation keeps a numeric-looking string from
# being coerced into a float (which would also silently drop values like
# the trailing zero in "21.70" or "0.006990").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRef, $value)
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# --- Price (column D) refresh ---
Set-TextValue "D2"  "235.61"
Set-TextValue "D3"  "21.70"
Set-TextValue "D4"  "5.351"
Set-TextValue "D5"  "0.05587"
Set-TextValue "D6"  "6.474"
Set-TextValue "D8"  "0.7995"
Set-TextValue "D9"  "1.039"
Set-TextValue "D10" "0.1393"
Set-TextValue "D11" "0.07314"
Set-TextValue "D12" "0.03185"
Set-TextValue "D13" "0.02955"
Set-TextValue "D14" "0.09253"
Set-TextValue "D15" "0.001663"
Set-TextValue "D17" "0.04792"
Set-TextValue "D18" "0.0005717"
Set-TextValue "D19" "0.006227"
Set-TextValue "D20" "0.005053"
Set-TextValue "D21" "0.001052"
Set-TextValue "D22" "0.0001503"
Set-TextValue "D24" "3.949"
Set-TextValue "D25" "2.201"
Set-TextValue "D40" "0.04107"
Set-TextValue "D41" "0.006990"

# --- Row 42 / 43: BKEXToken and CEJI swap ranking order ---
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D42" "0.003506"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D43" "0.1034"
$ws.Range("E43").Value = "42BKEXTokenBKK"

Set-TextValue "D44" "0.008810"
Set-TextValue "D45" "0.00005443"
Set-TextValue "D48" "0.03494"
Set-TextValue "D49" "0.00002104"
Set-TextValue "D50" "0.01012"

# --- Worst-in-24h label moves from BOLO (row 48) to One (row 18) ---
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E48").Value = "47BOLOBOLO"
